$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original layout (rows 20-24):
#   20: $GPGGA  | B..K = x
#   21: $GPVTG  | L,M = x
#   22: $KSXT   | B,C,E,G,H,J,L,M,N,O = x
#   23: $GPROT  | P = x
#   24: $GPTRA2 | B,G,H,K,M,N,O = x
#
# Target layout (rows 20-23), $GPTRA2 entry removed entirely:
#   20: $GPGGA  | B..K = x            (unchanged)
#   21: $KSXT   | L,M,N,O = x
#   22: $GPROT  | P = x
#   23: $GPVTG  | L,M = x

# Remove the old $KSXT row (row 22). This shifts the $GPROT row (23 -> 22)
# and the $GPTRA2 row (24 -> 23) up by one.
$ws.Rows(22).Delete()

# Row 21 currently holds the old $GPVTG entry (A21='$GPVTG', L21='x', M21='x').
# Re-purpose it as the new $KSXT entry, keeping L/M and adding N/O marks.
$ws.Range("A21").Value = '$KSXT'
$ws.Range("N21").Value = "x"
$ws.Range("O21").Value = "x"

# Row 22 now holds the shifted $GPROT entry (A22='$GPROT', P22='x') which
# already matches the target - no changes needed.

# Row 23 now holds the shifted $GPTRA2 entry. Delete it entirely and
# re-create it as the new $GPVTG entry with only L/M marked.
$ws.Rows(23).Delete()
$ws.Range("A23").Value = '$GPVTG'
$ws.Range("L23").Value = "x"
$ws.Range("M23").Value = "x"

# Select the (now empty) row below the table, matching the saved selection.
$ws.Rows(24).Select() | Out-Null
